$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new row's cells to be stored as text (matching the source data,
# which keeps trailing zeros like "4000.00"), then strip the formatting
# that forcing text-entry leaves behind so the cells end up with no
# explicit style, same as plain data rows in this sheet.
$row2 = $ws.Range("A2:G2")
$row2.NumberFormat = "@"

$ws.Range("A2").Value = "31/03/2001"
$ws.Range("B2").Value = "4000.00"
$ws.Range("C2").Value = "4000.00"
$ws.Range("D2").Value = "4000.00"
$ws.Range("E2").Value = "4000.00"
$ws.Range("F2").Value = "0.00"
$ws.Range("G2").Value = "100.00"

$row2.ClearFormats()
